$wb = $excel.ActiveWorkbook

# --- Fix the "#! FINISHED" -> "#! FINISH" marker text on the existing sheets ---
# (this also makes the now-unused "#! FINISHED" shared string get dropped)
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("C1").Value = "#! FINISH"
$ws2.Range("C2").Select() | Out-Null

$ws2b = $wb.Worksheets.Item("Sheet2 (2)")
$ws2b.Range("C1").Value = "#! FINISH"
$ws2b.Range("C2").Select() | Out-Null

$ws2c = $wb.Worksheets.Item("Sheet2 (3)")
$ws2c.Range("C1").Value = "#! FINISH"
$ws2c.Range("C2").Select() | Out-Null

$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Range("A1").Value = "#! FINISH"
$ws3.Range("A2").Select() | Out-Null

# --- Add the new worksheets used to test overly-long / missing property names ---

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSheet4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$wsSheet4.Name = "Sheet4"
$wsSheet4.Range("A1").Value = "#! WS_NAME doesntExistProperty"
$wsSheet4.Range("E1").Value = "#! FINISH"
$wsSheet4.Range("A1:XFD1").Select() | Out-Null

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSheet5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$wsSheet5.Name = "Sheet5"
$wsSheet5.Range("A1").Value = "#! WS_NAME doesntExistPropertyDuplicate"
$wsSheet5.Range("E1").Value = "#! FINISH"
$wsSheet5.Range("A1:XFD1").Select() | Out-Null

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSheet5b = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$wsSheet5b.Name = "Sheet5 (2)"
$wsSheet5b.Range("A1").Value = "#! WS_NAME doesntExistPropertyDuplicate"
$wsSheet5b.Range("E1").Value = "#! FINISH"
$wsSheet5b.Range("A1:XFD1").Select() | Out-Null

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSheet7 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$wsSheet7.Name = "Sheet7"
$wsSheet7.Range("A1").Value = "#! WS_NAME doesntExistPropertyButVeryLongAndShouldBeTruncated"
$wsSheet7.Range("I1").Value = "#! FINISH"
$wsSheet7.Range("I2").Select() | Out-Null

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSheet8 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$wsSheet8.Name = "Sheet8"
$wsSheet8.Range("A1").Value = "#! WS_NAME veryLong"
$wsSheet8.Range("F1").Value = "#! FINISH"
$wsSheet8.Range("A2").Select() | Out-Null
$wsSheet8.Activate() | Out-Null
